$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row (row 1) values B1:E1
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values B2:E2
$ws.Range("B2").Value = 15.847008530260693
$ws.Range("C2").Value = 11.244096314590454
$ws.Range("D2").Value = 16.214735093701755
$ws.Range("E2").Value = 9.4304240103828096

# Update row 3 values B3:E3
$ws.Range("B3").Value = 13.525050225365455
$ws.Range("C3").Value = 14.207823722985925
$ws.Range("D3").Value = 13.407086188382625
$ws.Range("E3").Value = 16.832541100581402

# Reflect the updated selection range as seen in the saved file
$ws.Range("B1:E3").Select()

$wb.Save()
